$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Rows.Item(3).Delete()

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Rows.Item(3).Delete()
$ws2.Range("E2").Value = "2016-03-22 12:55:43"
$ws2.Range("H2").Value = "2016-03-22 12:56:06"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Rows.Item(3).Delete()
$ws3.Range("E2").Value = "2016-03-22 12:55:47"
$ws3.Range("H2").Value = "2016-03-22 12:56:12"
